# The author reworked the opening/contribution paragraph of the draft after
# meeting with HP: it now leads with the NISAR mission description, explains
# the L-band penetration physics, and reorders the closing sentences. The
# paragraph spacing afterward was also tightened by one blank paragraph.

$d = $word.ActiveDocument

$newText = 'The forthcoming NASA-ISRO Synthetic Aperture Radar (NISAR) satellite mission presents a massive opportunity for snow monitoring. NISAR is equipped with L-band and S-band radar systems and will observe the Earth''s surface at a resolution of 3-10 meters (mode dependent), with a revisit frequency of twice every 12 days. Because the L-band radar operates at a low frequency, it can penetrate snow. As the wave moves through the snow, it slows down, causing a phase shift. This phase shift can be used to estimate changes in snow depth and Snow Water Equivalent. This work uses NISAR-like L-band products from the NASA JPL UAVSAR sensor to estimate total snow depth using machine learning (ML) algorithms, assuming that snow accumulation patterns are approximately consistent. Our findings validated our understanding of how the physics of L-band radar interacts with the snow. By showcasing the effectiveness of UAVSAR data in snow depth estimation, this work is a proof of concept for using L-band data to estimate total snow depth in preparation for the NISAR mission. Our findings are a step towards developing a global snow depth prediction system that provides valuable information for water resource management, flood forecasting, and avalanche hazard assessment, provided that accurate and representative training data is available.'

# Replace paragraph 1's text, but stop one character short of its end so we
# don't swallow the paragraph mark (which would merge it with paragraph 2).
$p1 = $d.Paragraphs(1).Range
$body = $d.Range($p1.Start, $p1.End - 1)
$body.Text = $newText

# Remove one of the blank paragraphs that used to follow the main paragraph.
$d.Paragraphs(2).Range.Delete()
